# Update the "想去人数" (interest count) figures in column F across the
# four worksheets to reflect the refreshed scrape, per the commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 71
$ws1.Range("F6").Value  = 372
$ws1.Range("F8").Value  = 561
$ws1.Range("F9").Value  = 1462
$ws1.Range("F11").Value = 1352
$ws1.Range("F12").Value = 3010
$ws1.Range("F13").Value = 438
$ws1.Range("F14").Value = 1633
$ws1.Range("F17").Value = 241
$ws1.Range("F18").Value = 1397
$ws1.Range("F19").Value = 268
$ws1.Range("F21").Value = 1129
$ws1.Range("F22").Value = 404
$ws1.Range("F23").Value = 5
$ws1.Range("F24").Value = 3503
$ws1.Range("F27").Value = 1549

# 演出 (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 170
$ws2.Range("F7").Value = 53
$ws2.Range("F8").Value = 23

# 本地生活 (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 797

# 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 797
$ws4.Range("F6").Value  = 71
$ws4.Range("F8").Value  = 170
$ws4.Range("F11").Value = 53
$ws4.Range("F13").Value = 23
$ws4.Range("F16").Value = 372
$ws4.Range("F18").Value = 561
$ws4.Range("F19").Value = 1462
$ws4.Range("F21").Value = 1352
$ws4.Range("F22").Value = 3010
$ws4.Range("F23").Value = 438
$ws4.Range("F24").Value = 1633
$ws4.Range("F27").Value = 241
$ws4.Range("F28").Value = 1397
$ws4.Range("F29").Value = 268
$ws4.Range("F33").Value = 1129
$ws4.Range("F34").Value = 404
$ws4.Range("F35").Value = 5
$ws4.Range("F36").Value = 3503
$ws4.Range("F39").Value = 1549
